$d = $word.ActiveDocument

# Locate the final (target) paragraph -- the last paragraph in the document,
# which currently holds only the _GoBack bookmark.
$count = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($count)
$startPos = $targetPara.Range.Start

# A collapsed range positioned exactly at the start of that paragraph; inserting
# WordOpenXML there adds new paragraphs ahead of it while the final inserted
# paragraph's runs merge into the existing (bookmarked) paragraph.
$insertionPoint = $d.Range($startPos, $startPos)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Note: the exact solution to this is rather difficult to derive.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Suppose 2 teams </w:t></w:r><w:r><w:t xml:space="preserve">A and B </w:t></w:r><w:r><w:t>are playing a series of games and the first team to win 4 games wins the series.</w:t></w:r><w:r><w:t xml:space="preserve">  Suppose that </w:t></w:r><w:r><w:t xml:space="preserve">team A has a </w:t></w:r><w:r><w:t>55% chance of winning each game and that the outcome of each game is independent.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>W</w:t></w:r><w:r><w:t>hat is the probability that team A wins the series?</w:t></w:r><w:r><w:t xml:space="preserve">  Give an exact result and confirm it via simulation.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>What</w:t></w:r><w:r><w:t xml:space="preserve"> is the </w:t></w:r><w:r><w:t>expected number of games played</w:t></w:r><w:r><w:t>?  Give an exact result and confirm it via simulation.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">What is the expected number of games played given that team A wins the series?  </w:t></w:r><w:r><w:t>Give an exact result and confirm it via simulation.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Now suppose we only know that team A is more likely to win each game, but do not know the exact probability.  If the most likely number of games played is 5, what does this imply about the probability that team A wins each game?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($xml)

# The merge above leaves the original paragraph's list level (ilvl 0) in place;
# the final question belongs one level deeper (ilvl 1), so fix it up.
$newCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($newCount)
$lastPara.Range.ListFormat.ListLevelNumber = 2
